$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price / volume data scraped on Sun May 19 13:49:54 UTC 2024

$ws.Range('D2').Value = '66.964.58'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '3.082.95'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.27'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.59'
$ws.Range('E6').Value = '  -2.99%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').Value = '3.081.24'
$ws.Range('E8').Value = '  -0.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.515'
$ws.Range('E9').Value = '  -1.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.42'
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.472'
$ws.Range('E12').Value = '  -1.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000241'
$ws.Range('E13').Value = '  -1.95%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.30'
$ws.Range('E14').Value = '  -2.26%  '
$ws.Range('E15').Value = '  -2.11%  '
$ws.Range('D16').Value = '3.597.78'
$ws.Range('E16').Value = '  -0.88%  '
$ws.Range('D17').Value = '66.840.58'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.01'
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('D19').Value = '3.078.58'
$ws.Range('E19').Value = '  -1.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.50'
$ws.Range('E20').Value = '  +1.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '488.97'
$ws.Range('E21').Value = '  +2.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.76'
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.687'
$ws.Range('E23').Value = '  -3.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.71'
$ws.Range('E24').Value = '  -1.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.88'
$ws.Range('E25').Value = '  -3.20%  '
$ws.Range('E26').Value = '  -2.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.30'
$ws.Range('E27').Value = '  +3.38%  '
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.80'
$ws.Range('E29').Value = '  -2.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.32'
$ws.Range('E30').Value = '  -3.95%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.63'
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.88'
$ws.Range('E32').Value = '  -2.60%  '
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('D34').Value = '0.0₃0912'
$ws.Range('E34').Value = '  -5.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.71'
$ws.Range('E36').Value = '  -2.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.955'
$ws.Range('E37').Value = '  -2.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '46.27'
$ws.Range('E38').Value = '  -3.13%  '
$ws.Range('E39').Value = '  +1.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.99'
$ws.Range('E40').Value = '  -4.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.303'
$ws.Range('E41').Value = '  -2.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.33'
$ws.Range('E42').Value = '  -3.06%  '
$ws.Range('D43').Value = '2.777.79'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '370.40'
$ws.Range('E44').Value = '  -2.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0346'
$ws.Range('E45').Value = '  -2.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '135.62'
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.48'
$ws.Range('E47').Value = '  -3.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '24.43'
$ws.Range('E49').Value = '  -1.32%  '
$ws.Range('E50').Value = '  -2.22%  '
$ws.Range('E51').Value = '  -1.24%  '
